$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Swap the data (columns B:AC) between pairs of rows that were re-sorted -----------
# (column A holds the running row index and must stay put)
function Swap-RowData($sheet, $row1, $row2) {
    $range1 = $sheet.Range("B$row1`:AC$row1")
    $range2 = $sheet.Range("B$row2`:AC$row2")
    $tmp = $range1.Value2
    $range1.Value2 = $range2.Value2
    $range2.Value2 = $tmp
}

Swap-RowData $ws 91 92
Swap-RowData $ws 110 111
Swap-RowData $ws 129 131

# --- Append the new match record as row 194 --------------------------------------------
# Clone the formatting of the last existing row (193) onto the new row first.
$ws.Range("A193:AC193").Copy() | Out-Null
$ws.Range("A194:AC194").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A194").Value = 192
$ws.Range("B194").Value = 7623986
$ws.Range("C194").Value = "Costa Rica Primera Division"
$ws.Range("D194").Value = "Costa Rica Primera Division"
$ws.Range("E194").Value = 45347.875
$ws.Range("F194").Value = "Municipal Perez Zeledon"
$ws.Range("G194").Value = "Cartagines"
$ws.Range("H194").Value = 2
$ws.Range("I194").Value = 2
$ws.Range("J194").Value = "D"
$ws.Range("K194").Value = 2.625
$ws.Range("L194").Value = 3.3
$ws.Range("M194").Value = 2.375
$ws.Range("N194").Value = 3
$ws.Range("O194").Value = 3.3
$ws.Range("P194").Value = 2.15
$ws.Range("Q194").Value = 0.25
$ws.Range("R194").Value = 1.85
$ws.Range("S194").Value = 1.95
$ws.Range("T194").Value = 2.5
$ws.Range("U194").Value = 1.85
$ws.Range("V194").Value = 1.95
$ws.Range("W194").Value = -1
$ws.Range("X194").Value = 2.3
$ws.Range("Y194").Value = -1
$ws.Range("Z194").Value = 0.425
$ws.Range("AA194").Value = -0.5
$ws.Range("AB194").Value = 0.8500000000000001
$ws.Range("AC194").Value = -1
